# Update countries & provincias Spain
#
# - Refresh case/death figures for Hungria, Vietnam and the
#   Nueva Zelanda / Uzbekistan pair.
# - Re-seat "Suazilandia" (Eswatini) with fresh figures just above
#   "Timor Oriental", shifting Timor Oriental / Belice / Nueva Caledonia /
#   Islas Virgenes de los Estados Unidos / Fiyi / Malaui down one row,
#   removing the old "Suazilandia" row that used to sit below Malaui.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 62: Hungria
$ws.Range("A62").Value = 'Hungria'
$ws.Range("B62").Value = 1834
$ws.Range("C62").Value = 71
$ws.Range("D62").Value = 231
$ws.Range("E62").Value = 1431
$ws.Range("F62").Value = 60
$ws.Range("G62").Value = 16
$ws.Range("H62").Value = 172

# Row 70: Uzbekistan
$ws.Range("A70").Value = 'Uzbekistan'
$ws.Range("B70").Value = 1450
$ws.Range("C70").Value = 45
$ws.Range("D70").Value = 156
$ws.Range("E70").Value = 1290
$ws.Range("F70").Value = 8
$ws.Range("G70").Value = 0
$ws.Range("H70").Value = 4

# Row 71: Nueva Zelanda
$ws.Range("A71").Value = 'Nueva Zelanda'
$ws.Range("B71").Value = 1422
$ws.Range("C71").Value = 13
$ws.Range("D71").Value = 867
$ws.Range("E71").Value = 544
$ws.Range("F71").Value = 3
$ws.Range("G71").Value = 0
$ws.Range("H71").Value = 11

# Row 116: Vietnam
$ws.Range("A116").Value = 'Vietnam'
$ws.Range("B116").Value = 268
$ws.Range("C116").Value = 0
$ws.Range("D116").Value = 201
$ws.Range("E116").Value = 67
$ws.Range("F116").Value = 8
$ws.Range("G116").Value = 0
$ws.Range("H116").Value = 0

# Row 179: Suazilandia
$ws.Range("A179").Value = 'Suazilandia'
$ws.Range("B179").Value = 19
$ws.Range("C179").Value = 3
$ws.Range("D179").Value = 8
$ws.Range("E179").Value = 10
$ws.Range("F179").Value = 0
$ws.Range("G179").Value = 0
$ws.Range("H179").Value = 1

# Row 180: Timor Oriental
$ws.Range("A180").Value = 'Timor Oriental'
$ws.Range("B180").Value = 18
$ws.Range("C180").Value = 0
$ws.Range("D180").Value = 1
$ws.Range("E180").Value = 17
$ws.Range("F180").Value = 0
$ws.Range("G180").Value = 0
$ws.Range("H180").Value = 0

# Row 181: Belice
$ws.Range("A181").Value = 'Belice'
$ws.Range("B181").Value = 18
$ws.Range("C181").Value = 0
$ws.Range("D181").Value = 0
$ws.Range("E181").Value = 16
$ws.Range("F181").Value = 1
$ws.Range("G181").Value = 0
$ws.Range("H181").Value = 2

# Row 182: Nueva Caledonia
$ws.Range("A182").Value = 'Nueva Caledonia'
$ws.Range("B182").Value = 18
$ws.Range("C182").Value = 0
$ws.Range("D182").Value = 14
$ws.Range("E182").Value = 4
$ws.Range("F182").Value = 1
$ws.Range("G182").Value = 0
$ws.Range("H182").Value = 0

# Row 183: Islas Virgenes de los Estados Unidos
$ws.Range("A183").Value = 'Islas Virgenes de los Estados Unidos'
$ws.Range("B183").Value = 17
$ws.Range("C183").Value = 0
$ws.Range("D183").Value = 0
$ws.Range("E183").Value = 17
$ws.Range("F183").Value = 0
$ws.Range("G183").Value = 0
$ws.Range("H183").Value = 0

# Row 184: Fiyi
$ws.Range("A184").Value = 'Fiyi'
$ws.Range("B184").Value = 17
$ws.Range("C184").Value = 0
$ws.Range("D184").Value = 0
$ws.Range("E184").Value = 17
$ws.Range("F184").Value = 0
$ws.Range("G184").Value = 0
$ws.Range("H184").Value = 0

# Row 185: Malaui
$ws.Range("A185").Value = 'Malaui'
$ws.Range("B185").Value = 17
$ws.Range("C185").Value = 0
$ws.Range("D185").Value = 3
$ws.Range("E185").Value = 12
$ws.Range("F185").Value = 1
$ws.Range("G185").Value = 0
$ws.Range("H185").Value = 2

# Row 186: Namibia
$ws.Range("A186").Value = 'Namibia'
$ws.Range("B186").Value = 16
$ws.Range("C186").Value = 0
$ws.Range("D186").Value = 4
$ws.Range("E186").Value = 12
$ws.Range("F186").Value = 0
$ws.Range("G186").Value = 0
$ws.Range("H186").Value = 0
